# Natmi following Dr Hou advice
# Ligand-expressing-cells / receptor-expressing-cells counts (columns E and K)
# go from 1 to 3 for every data row, and the dependent expression /
# specificity / edge-weight columns are updated to the recomputed values
# that follow from that change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.402677
$ws.Range("H2").Value = 37.20803100000001
$ws.Range("I2").Value = 0.4952943482020729
$ws.Range("J2").Value = 0.495294348202073
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.840326
$ws.Range("N2").Value = 38.520978
$ws.Range("O2").Value = 0.3393128690704512
$ws.Range("P2").Value = 0.3393128690704511
$ws.Range("Q2").Value = 159.254415952702
$ws.Range("R2").Value = 1433.289743574318
$ws.Range("S2").Value = 0.1680597463228244
$ws.Range("T2").Value = 0.1680597463228244

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.402677
$ws.Range("H3").Value = 37.20803100000001
$ws.Range("I3").Value = 0.4952943482020729
$ws.Range("J3").Value = 0.495294348202073
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.834223333333334
$ws.Range("N3").Value = 29.50267
$ws.Range("O3").Value = 0.2598749077175229
$ws.Range("P3").Value = 0.2598749077175228
$ws.Range("Q3").Value = 121.9706955491967
$ws.Range("R3").Value = 1097.73625994277
$ws.Range("S3").Value = 0.1287145730320243
$ws.Range("T3").Value = 0.1287145730320243

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.402677
$ws.Range("H4").Value = 37.20803100000001
$ws.Range("I4").Value = 0.4952943482020729
$ws.Range("J4").Value = 0.495294348202073
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.654269333333334
$ws.Range("N4").Value = 25.962808
$ws.Range("O4").Value = 0.228693956617749
$ws.Range("P4").Value = 0.2286939566177489
$ws.Range("Q4").Value = 107.3361072123387
$ws.Range("R4").Value = 966.0249649110481
$ws.Range("S4").Value = 0.1132708241807411
$ws.Range("T4").Value = 0.1132708241807411

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.402677
$ws.Range("H5").Value = 37.20803100000001
$ws.Range("I5").Value = 0.4952943482020729
$ws.Range("J5").Value = 0.495294348202073
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.513324000000001
$ws.Range("N5").Value = 19.539972
$ws.Range("O5").Value = 0.1721182665942771
$ws.Range("P5").Value = 0.1721182665942771
$ws.Range("Q5").Value = 80.78265376834803
$ws.Range("R5").Value = 727.0438839151321
$ws.Range("S5").Value = 0.0852492046664831
$ws.Range("T5").Value = 0.08524920466648309

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.161818666666666
$ws.Range("H6").Value = 15.485456
$ws.Range("I6").Value = 0.2061344991927113
$ws.Range("J6").Value = 0.2061344991927113
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 12.840326
$ws.Range("N6").Value = 38.520978
$ws.Range("O6").Value = 0.3393128690704512
$ws.Range("P6").Value = 0.3393128690704511
$ws.Range("Q6").Value = 66.27943443288532
$ws.Range("R6").Value = 596.514909895968
$ws.Range("S6").Value = 0.06994408833547948
$ws.Range("T6").Value = 0.06994408833547948

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.161818666666666
$ws.Range("H7").Value = 15.485456
$ws.Range("I7").Value = 0.2061344991927113
$ws.Range("J7").Value = 0.2061344991927113
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.834223333333334
$ws.Range("N7").Value = 29.50267
$ws.Range("O7").Value = 0.2598749077175229
$ws.Range("P7").Value = 0.2598749077175228
$ws.Range("Q7").Value = 50.76247757416888
$ws.Range("R7").Value = 456.86229816752
$ws.Range("S7").Value = 0.05356918395510363
$ws.Range("T7").Value = 0.05356918395510364

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.161818666666666
$ws.Range("H8").Value = 15.485456
$ws.Range("I8").Value = 0.2061344991927113
$ws.Range("J8").Value = 0.2061344991927113
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.654269333333334
$ws.Range("N8").Value = 25.962808
$ws.Range("O8").Value = 0.228693956617749
$ws.Range("P8").Value = 0.2286939566177489
$ws.Range("Q8").Value = 44.67176899116089
$ws.Range("R8").Value = 402.0459209204479
$ws.Range("S8").Value = 0.04714171421579933
$ws.Range("T8").Value = 0.04714171421579932

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.161818666666666
$ws.Range("H9").Value = 15.485456
$ws.Range("I9").Value = 0.2061344991927113
$ws.Range("J9").Value = 0.2061344991927113
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.513324000000001
$ws.Range("N9").Value = 19.539972
$ws.Range("O9").Value = 0.1721182665942771
$ws.Range("P9").Value = 0.1721182665942771
$ws.Range("Q9").Value = 33.620597405248
$ws.Range("R9").Value = 302.585376647232
$ws.Range("S9").Value = 0.03547951268632888
$ws.Range("T9").Value = 0.03547951268632888

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.073119
$ws.Range("H10").Value = 9.219357
$ws.Range("I10").Value = 0.1227233823836907
$ws.Range("J10").Value = 0.1227233823836907
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.840326
$ws.Range("N10").Value = 38.520978
$ws.Range("O10").Value = 0.3393128690704512
$ws.Range("P10").Value = 0.3393128690704511
$ws.Range("Q10").Value = 39.459849796794
$ws.Range("R10").Value = 355.138648171146
$ws.Range("S10").Value = 0.04164162297864017
$ws.Range("T10").Value = 0.04164162297864016

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.073119
$ws.Range("H11").Value = 9.219357
$ws.Range("I11").Value = 0.1227233823836907
$ws.Range("J11").Value = 0.1227233823836907
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.834223333333334
$ws.Range("N11").Value = 29.50267
$ws.Range("O11").Value = 0.2598749077175229
$ws.Range("P11").Value = 0.2598749077175228
$ws.Range("Q11").Value = 30.22173857591
$ws.Range("R11").Value = 271.99564718319
$ws.Range("S11").Value = 0.03189272767174389
$ws.Range("T11").Value = 0.03189272767174389

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.073119
$ws.Range("H12").Value = 9.219357
$ws.Range("I12").Value = 0.1227233823836907
$ws.Range("J12").Value = 0.1227233823836907
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 8.654269333333334
$ws.Range("N12").Value = 25.962808
$ws.Range("O12").Value = 0.228693956617749
$ws.Range("P12").Value = 0.2286939566177489
$ws.Range("Q12").Value = 26.595599519384
$ws.Range("R12").Value = 239.360395674456
$ws.Range("S12").Value = 0.02806609588683918
$ws.Range("T12").Value = 0.02806609588683918

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.073119
$ws.Range("H13").Value = 9.219357
$ws.Range("I13").Value = 0.1227233823836907
$ws.Range("J13").Value = 0.1227233823836907
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.513324000000001
$ws.Range("N13").Value = 19.539972
$ws.Range("O13").Value = 0.1721182665942771
$ws.Range("P13").Value = 0.1721182665942771
$ws.Range("Q13").Value = 20.016219737556
$ws.Range("R13").Value = 180.145977638004
$ws.Range("S13").Value = 0.02112293584646749
$ws.Range("T13").Value = 0.02112293584646749

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.403408
$ws.Range("H14").Value = 13.210224
$ws.Range("I14").Value = 0.175847770221525
$ws.Range("J14").Value = 0.175847770221525
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 12.840326
$ws.Range("N14").Value = 38.520978
$ws.Range("O14").Value = 0.3393128690704512
$ws.Range("P14").Value = 0.3393128690704511
$ws.Range("Q14").Value = 56.54119423100799
$ws.Range("R14").Value = 508.870748079072
$ws.Range("S14").Value = 0.0596674114335071
$ws.Range("T14").Value = 0.05966741143350709

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.403408
$ws.Range("H15").Value = 13.210224
$ws.Range("I15").Value = 0.175847770221525
$ws.Range("J15").Value = 0.175847770221525
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.834223333333334
$ws.Range("N15").Value = 29.50267
$ws.Range("O15").Value = 0.2598749077175229
$ws.Range("P15").Value = 0.2598749077175228
$ws.Range("Q15").Value = 43.30409769978667
$ws.Range("R15").Value = 389.7368792980801
$ws.Range("S15").Value = 0.04569842305865097
$ws.Range("T15").Value = 0.04569842305865097

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.403408
$ws.Range("H16").Value = 13.210224
$ws.Range("I16").Value = 0.175847770221525
$ws.Range("J16").Value = 0.175847770221525
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 8.654269333333334
$ws.Range("N16").Value = 25.962808
$ws.Range("O16").Value = 0.228693956617749
$ws.Range("P16").Value = 0.2286939566177489
$ws.Range("Q16").Value = 38.10827881655467
$ws.Range("R16").Value = 342.974509348992
$ws.Range("S16").Value = 0.04021532233436934
$ws.Range("T16").Value = 0.04021532233436933

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.403408
$ws.Range("H17").Value = 13.210224
$ws.Range("I17").Value = 0.175847770221525
$ws.Range("J17").Value = 0.175847770221525
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.513324000000001
$ws.Range("N17").Value = 19.539972
$ws.Range("O17").Value = 0.1721182665942771
$ws.Range("P17").Value = 0.1721182665942771
$ws.Range("Q17").Value = 28.680823008192
$ws.Range("R17").Value = 258.127407073728
$ws.Range("S17").Value = 0.03026661339499763
$ws.Range("T17").Value = 0.03026661339499762
